# Insert 4 new rows above the existing row 362, shifting all following
# rows down by 4 (this preserves their content/formatting automatically).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("362:365").Insert()

# Populate the 4 newly-inserted rows with the new weekly price records.
# Row 362: Naranja - Lane Late - Primera
$ws.Cells.Item(362, 1).Value  = 3
$ws.Cells.Item(362, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(362, 3).Value  = "Coquimbo"
$ws.Cells.Item(362, 4).Value  = 44449
$ws.Cells.Item(362, 5).Value  = 5
$ws.Cells.Item(362, 6).Value  = "Fruta"
$ws.Cells.Item(362, 7).Value  = 100102
$ws.Cells.Item(362, 8).Value  = "Cítricos"
$ws.Cells.Item(362, 9).Value  = 100102005
$ws.Cells.Item(362, 10).Value = "Naranja"
$ws.Cells.Item(362, 11).Value = "Lane Late"
$ws.Cells.Item(362, 12).Value = "Primera"
$ws.Cells.Item(362, 13).Value = 170
$ws.Cells.Item(362, 14).Value = 3500
$ws.Cells.Item(362, 15).Value = 4000
$ws.Cells.Item(362, 16).Value = 3735
$ws.Cells.Item(362, 17).Value = "$/malla 13 kilos"
$ws.Cells.Item(362, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(362, 19).Value = 287
$ws.Cells.Item(362, 20).Value = 13

# Row 363: Naranja - Lane Late - Segunda
$ws.Cells.Item(363, 1).Value  = 3
$ws.Cells.Item(363, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(363, 3).Value  = "Coquimbo"
$ws.Cells.Item(363, 4).Value  = 44449
$ws.Cells.Item(363, 5).Value  = 5
$ws.Cells.Item(363, 6).Value  = "Fruta"
$ws.Cells.Item(363, 7).Value  = 100102
$ws.Cells.Item(363, 8).Value  = "Cítricos"
$ws.Cells.Item(363, 9).Value  = 100102005
$ws.Cells.Item(363, 10).Value = "Naranja"
$ws.Cells.Item(363, 11).Value = "Lane Late"
$ws.Cells.Item(363, 12).Value = "Segunda"
$ws.Cells.Item(363, 13).Value = 85
$ws.Cells.Item(363, 14).Value = 3000
$ws.Cells.Item(363, 15).Value = 3000
$ws.Cells.Item(363, 16).Value = 3000
$ws.Cells.Item(363, 17).Value = "$/malla 13 kilos"
$ws.Cells.Item(363, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(363, 19).Value = 231
$ws.Cells.Item(363, 20).Value = 13

# Row 364: Naranja - Navel Late - Primera
$ws.Cells.Item(364, 1).Value  = 3
$ws.Cells.Item(364, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(364, 3).Value  = "Coquimbo"
$ws.Cells.Item(364, 4).Value  = 44449
$ws.Cells.Item(364, 5).Value  = 5
$ws.Cells.Item(364, 6).Value  = "Fruta"
$ws.Cells.Item(364, 7).Value  = 100102
$ws.Cells.Item(364, 8).Value  = "Cítricos"
$ws.Cells.Item(364, 9).Value  = 100102005
$ws.Cells.Item(364, 10).Value = "Naranja"
$ws.Cells.Item(364, 11).Value = "Navel Late"
$ws.Cells.Item(364, 12).Value = "Primera"
$ws.Cells.Item(364, 13).Value = 155
$ws.Cells.Item(364, 14).Value = 3500
$ws.Cells.Item(364, 15).Value = 4000
$ws.Cells.Item(364, 16).Value = 3758
$ws.Cells.Item(364, 17).Value = "$/malla 13 kilos"
$ws.Cells.Item(364, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(364, 19).Value = 289
$ws.Cells.Item(364, 20).Value = 13

# Row 365: Naranja - Navel Late - Segunda
$ws.Cells.Item(365, 1).Value  = 3
$ws.Cells.Item(365, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(365, 3).Value  = "Coquimbo"
$ws.Cells.Item(365, 4).Value  = 44449
$ws.Cells.Item(365, 5).Value  = 5
$ws.Cells.Item(365, 6).Value  = "Fruta"
$ws.Cells.Item(365, 7).Value  = 100102
$ws.Cells.Item(365, 8).Value  = "Cítricos"
$ws.Cells.Item(365, 9).Value  = 100102005
$ws.Cells.Item(365, 10).Value = "Naranja"
$ws.Cells.Item(365, 11).Value = "Navel Late"
$ws.Cells.Item(365, 12).Value = "Segunda"
$ws.Cells.Item(365, 13).Value = 80
$ws.Cells.Item(365, 14).Value = 3000
$ws.Cells.Item(365, 15).Value = 3000
$ws.Cells.Item(365, 16).Value = 3000
$ws.Cells.Item(365, 17).Value = "$/malla 13 kilos"
$ws.Cells.Item(365, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(365, 19).Value = 231
$ws.Cells.Item(365, 20).Value = 13
